$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 16.149999989999998
$ws.Range("C2").Value = 13.249999989999999
$ws.Range("D2").Value = 16.149999989999998
$ws.Range("E2").Value = 13.249999989999999

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 11.949999989999998
$ws.Range("C3").Value = 6.8499999899999997
$ws.Range("D3").Value = 11.949999989999998
$ws.Range("E3").Value = 6.8499999899999997

# Update the selected range to match the new relevant data extent
$ws.Range("B1:E3").Select() | Out-Null
